# Updated Excel & credential
# - Adds two new TestCases rows (ToggleTc005, ToggleTc006)
# - Marks the previous last row ("No" instead of "Yes")
# - Extends conditional formatting + data validation to the new rows

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# --- 1. Existing row 10 (ToggleTc003_1_CreatePreconditionData) flips C10 to "No" ---
$ws.Range("C10").Value2 = "No"

# --- 2. New row 11: ToggleTc005_... ---
$ws.Range("A11").Value2 = "ToggleTc005_VerifyTheToggleButtonAndFunctionalityOfToggleButtonOnInstitutionPageRecord"
$ws.Range("B11").Value2 = 1
$ws.Range("C11").Value2 = "No"
$ws.Range("D11").Value2 = "Skip: Disabled in excel."
$ws.Range("E11").Value2 = "High"

# --- 3. New row 12: ToggleTc006_... ---
$ws.Range("A12").Value2 = "ToggleTc006_CreateNewCustomSDG"
$ws.Range("B12").Value2 = 1
$ws.Range("C12").Value2 = "No"
$ws.Range("D12").Value2 = "Skip: Disabled in excel."
$ws.Range("E12").Value2 = "High"

# --- 4. Copy the formatting (styles/borders/fonts/number formats) of row 10 onto
#        the two new rows so they look identical to the existing data rows. ---
$ws.Range("A10:E10").Copy()
$ws.Range("A11:E11").PasteSpecial(-4122)
$ws.Range("A10:E10").Copy()
$ws.Range("A12:E12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 5. Replicate the conditional formatting rules used on every data row onto
#        the new rows 11 and 12 (E column: High/Low priority; D column: Pass/
#        Fail/Skip banners, duplicated exactly as the workbook already does for
#        every other row). Existing rules occupy priorities 1-96, so the new
#        ones start right after that (kept unique on purpose). ---
$nextPriority = 97

foreach ($row in 11, 12) {
  $eRange = $ws.Range("E$row")
  $eFc = $eRange.FormatConditions

  $cLow = $eFc.Add(1, 3, '"LOW"')
  $cLow.Font.Color = 24832
  $cLow.Interior.Color = 13561798
  $cLow.Priority = $nextPriority
  $nextPriority = $nextPriority + 1

  $cHigh = $eFc.Add(1, 3, '"High"')
  $cHigh.Font.Color = 393372
  $cHigh.Interior.Color = 13551615
  $cHigh.Priority = $nextPriority
  $nextPriority = $nextPriority + 1

  $cLow2 = $eFc.Add(1, 3, '"Low"')
  $cLow2.Interior.Color = 3506772
  $cLow2.Priority = $nextPriority
  $nextPriority = $nextPriority + 1

  $dRange = $ws.Range("D$row")

  # Two identical blocks, matching the pattern already present for every
  # other row in this sheet (Skip:/Fail/Pass banners).
  for ($block = 1; $block -le 2; $block++) {
    $dFc = $dRange.FormatConditions

    $cSkip = $dFc.Add(9, 0, [System.Type]::Missing, [System.Type]::Missing, "Skip:", 0)
    $cSkip.Font.Bold = $true
    $cSkip.Font.Color = 0
    $cSkip.Interior.Color = 16247774
    $cSkip.Priority = $nextPriority
    $nextPriority = $nextPriority + 1

    $cFail = $dFc.Add(9, 0, [System.Type]::Missing, [System.Type]::Missing, "Fail", 0)
    $cFail.Font.Bold = $true
    $cFail.Font.Color = 0
    $cFail.Interior.Color = 11389944
    $cFail.Priority = $nextPriority
    $nextPriority = $nextPriority + 1

    $cPass = $dFc.Add(9, 0, [System.Type]::Missing, [System.Type]::Missing, "Pass", 0)
    $cPass.Font.Bold = $true
    $cPass.Font.Color = 0
    $cPass.Interior.Color = 11854021
    $cPass.Priority = $nextPriority
    $nextPriority = $nextPriority + 1
  }
}

# --- 6. Extend the two data-validation list rules (Yes/No on column C, High/Low
#        on column E) to cover the two new rows. ---
$ws.Range("C2:C10").Validation.Delete()
$cVal = $ws.Range("C2:C12").Validation
$cVal.Add(3, 1, 1, '"Yes,No"')
$cVal.IgnoreBlank = $false

$ws.Range("E2:E10").Validation.Delete()
$eVal = $ws.Range("E2:E12").Validation
$eVal.Add(3, 1, 1, '"High,Low"')

# --- 7. Move the active selection the way the original author left it. ---
$ws.Range("C11").Select()
